$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2556606666666667
$ws.Range("H2").Value = 0.7669820000000001
$ws.Range("I2").Value = 0.01354513404628681
$ws.Range("J2").Value = 0.01354513404628681
$ws.Range("M2").Value = 51.15371566666666
$ws.Range("N2").Value = 153.461147
$ws.Range("O2").Value = 0.3311207986511828
$ws.Range("P2").Value = 0.3311207986511828
$ws.Range("Q2").Value = 13.07799304981711
$ws.Range("R2").Value = 117.701937448354
$ws.Range("S2").Value = 0.004485075603243815
$ws.Range("T2").Value = 0.004485075603243815
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2556606666666667
$ws.Range("H3").Value = 0.7669820000000001
$ws.Range("I3").Value = 0.01354513404628681
$ws.Range("J3").Value = 0.01354513404628681
$ws.Range("M3").Value = 53.36146666666667
$ws.Range("O3").Value = 0.3454116915964105
$ws.Range("P3").Value = 0.3454116915964106
$ws.Range("Q3").Value = 13.64242814231111
$ws.Range("R3").Value = 122.7818532808
$ws.Range("S3").Value = 0.004678647663828058
$ws.Range("T3").Value = 0.004678647663828059
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2556606666666667
$ws.Range("H4").Value = 0.7669820000000001
$ws.Range("I4").Value = 0.01354513404628681
$ws.Range("J4").Value = 0.01354513404628681
$ws.Range("M4").Value = 25.00653133333333
$ws.Range("N4").Value = 75.019594
$ws.Range("O4").Value = 0.1618686447050176
$ws.Range("P4").Value = 0.1618686447050176
$ws.Range("Q4").Value = 6.393186471700889
$ws.Range("R4").Value = 57.538678245308
$ws.Range("S4").Value = 0.002192532490420237
$ws.Range("T4").Value = 0.002192532490420237
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2556606666666667
$ws.Range("H5").Value = 0.7669820000000001
$ws.Range("I5").Value = 0.01354513404628681
$ws.Range("J5").Value = 0.01354513404628681
$ws.Range("M5").Value = 6.481347
$ws.Range("N5").Value = 19.444041
$ws.Range("O5").Value = 0.0419541135381084
$ws.Range("P5").Value = 0.0419541135381084
$ws.Range("Q5").Value = 1.657025494918
$ws.Range("R5").Value = 14.913229454262
$ws.Range("S5").Value = 0.0005682740916668143
$ws.Range("T5").Value = 0.0005682740916668143
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2556606666666667
$ws.Range("H6").Value = 0.7669820000000001
$ws.Range("I6").Value = 0.01354513404628681
$ws.Range("J6").Value = 0.01354513404628681
$ws.Range("M6").Value = 18.483507
$ws.Range("N6").Value = 55.450521
$ws.Range("O6").Value = 0.1196447515092806
$ws.Range("P6").Value = 0.1196447515092806
$ws.Range("Q6").Value = 4.725505721958
$ws.Range("R6").Value = 42.529551497622
$ws.Range("S6").Value = 0.001620604197127882
$ws.Range("T6").Value = 0.001620604197127882
$ws.Range("I7").Value = 0.8835639662863414
$ws.Range("J7").Value = 0.8835639662863415
$ws.Range("M7").Value = 51.15371566666666
$ws.Range("N7").Value = 153.461147
$ws.Range("O7").Value = 0.3311207986511828
$ws.Range("P7").Value = 0.3311207986511828
$ws.Range("Q7").Value = 853.0918461696072
$ws.Range("R7").Value = 7677.826615526465
$ws.Range("S7").Value = 0.2925664061761401
$ws.Range("T7").Value = 0.2925664061761402
$ws.Range("I8").Value = 0.8835639662863414
$ws.Range("J8").Value = 0.8835639662863415
$ws.Range("M8").Value = 53.36146666666667
$ws.Range("O8").Value = 0.3454116915964105
$ws.Range("P8").Value = 0.3454116915964106
$ws.Range("Q8").Value = 889.9105669981333
$ws.Range("R8").Value = 8009.1951029832
$ws.Range("S8").Value = 0.305193324228599
$ws.Range("T8").Value = 0.3051933242285991
$ws.Range("I9").Value = 0.8835639662863414
$ws.Range("J9").Value = 0.8835639662863415
$ws.Range("M9").Value = 25.00653133333333
$ws.Range("N9").Value = 75.019594
$ws.Range("O9").Value = 0.1618686447050176
$ws.Range("P9").Value = 0.1618686447050176
$ws.Range("Q9").Value = 417.0345732158146
$ws.Range("R9").Value = 3753.311158942332
$ws.Range("S9").Value = 0.14302130173296
$ws.Range("T9").Value = 0.14302130173296
$ws.Range("I10").Value = 0.8835639662863414
$ws.Range("J10").Value = 0.8835639662863415
$ws.Range("M10").Value = 6.481347
$ws.Range("N10").Value = 19.444041
$ws.Range("O10").Value = 0.0419541135381084
$ws.Range("P10").Value = 0.0419541135381084
$ws.Range("Q10").Value = 108.089592434022
$ws.Range("R10").Value = 972.8063319061978
$ws.Range("S10").Value = 0.03706914295975855
$ws.Range("T10").Value = 0.03706914295975855
$ws.Range("I11").Value = 0.8835639662863414
$ws.Range("J11").Value = 0.8835639662863415
$ws.Range("M11").Value = 18.483507
$ws.Range("N11").Value = 55.450521
$ws.Range("O11").Value = 0.1196447515092806
$ws.Range("P11").Value = 0.1196447515092806
$ws.Range("Q11").Value = 308.2499268101819
$ws.Range("R11").Value = 2774.249341291638
$ws.Range("S11").Value = 0.1057137911888837
$ws.Range("T11").Value = 0.1057137911888837
$ws.Range("G12").Value = 1.942037333333333
$ws.Range("H12").Value = 5.826112
$ws.Range("I12").Value = 0.1028908996673717
$ws.Range("J12").Value = 0.1028908996673717
$ws.Range("M12").Value = 51.15371566666666
$ws.Range("N12").Value = 153.461147
$ws.Range("O12").Value = 0.3311207986511828
$ws.Range("P12").Value = 0.3311207986511828
$ws.Range("Q12").Value = 99.34242556338488
$ws.Range("R12").Value = 894.081830070464
$ws.Range("S12").Value = 0.03406931687179885
$ws.Range("T12").Value = 0.03406931687179886
$ws.Range("G13").Value = 1.942037333333333
$ws.Range("H13").Value = 5.826112
$ws.Range("I13").Value = 0.1028908996673717
$ws.Range("J13").Value = 0.1028908996673717
$ws.Range("M13").Value = 53.36146666666667
$ws.Range("O13").Value = 0.3454116915964105
$ws.Range("P13").Value = 0.3454116915964106
$ws.Range("Q13").Value = 103.6299604280889
$ws.Range("R13").Value = 932.6696438528002
$ws.Range("S13").Value = 0.03553971970398342
$ws.Range("T13").Value = 0.03553971970398343
$ws.Range("G14").Value = 1.942037333333333
$ws.Range("H14").Value = 5.826112
$ws.Range("I14").Value = 0.1028908996673717
$ws.Range("J14").Value = 0.1028908996673717
$ws.Range("M14").Value = 25.00653133333333
$ws.Range("N14").Value = 75.019594
$ws.Range("O14").Value = 0.1618686447050176
$ws.Range("P14").Value = 0.1618686447050176
$ws.Range("Q14").Value = 48.56361742650311
$ws.Range("R14").Value = 437.072556838528
$ws.Range("S14").Value = 0.01665481048163741
$ws.Range("T14").Value = 0.01665481048163741
$ws.Range("G15").Value = 1.942037333333333
$ws.Range("H15").Value = 5.826112
$ws.Range("I15").Value = 0.1028908996673717
$ws.Range("J15").Value = 0.1028908996673717
$ws.Range("M15").Value = 6.481347
$ws.Range("N15").Value = 19.444041
$ws.Range("O15").Value = 0.0419541135381084
$ws.Range("P15").Value = 0.0419541135381084
$ws.Range("Q15").Value = 12.587017844288
$ws.Range("R15").Value = 113.283160598592
$ws.Range("S15").Value = 0.004316696486683033
$ws.Range("T15").Value = 0.004316696486683034
$ws.Range("G16").Value = 1.942037333333333
$ws.Range("H16").Value = 5.826112
$ws.Range("I16").Value = 0.1028908996673717
$ws.Range("J16").Value = 0.1028908996673717
$ws.Range("M16").Value = 18.483507
$ws.Range("N16").Value = 55.450521
$ws.Range("O16").Value = 0.1196447515092806
$ws.Range("P16").Value = 0.1196447515092806
$ws.Range("Q16").Value = 35.895660644928
$ws.Range("R16").Value = 323.060945804352
$ws.Range("S16").Value = 0.01231035612326901
$ws.Range("T16").Value = 0.01231035612326902
